# Add a 'hole_id' index column (A) to the 'train' sheet so that
# cross validation can line rows back up with their source hole.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
    "MHZ_12_04",
    "LBU_02_3",
    "MHZ_08_01",
    "LBU_05_02",
    "LBU_05_17",
    "LBU_01_3",
    "LBU_05_15",
    "LBU_05_29",
    "LBU_05_19",
    "LBU_05_24",
    "LBU_02_4",
    "LBU_05_25",
    "LBU_01_1",
    "LBU_05_22",
    "LBU_05_01",
    "MHZ_12_03",
    "LBU_87_5",
    "LBU_96_3",
    "LBU_05_30",
    "MHZ_12_01",
    "MHZ_08_05",
    "LBU_05_14",
    "LBU_05_18",
    "LBU_05_16",
    "LBU_98_2",
    "LBU_01_2",
    "LBU_05_07",
    "LBU_05_09",
    "LBU_98_6",
    "MHZ_12_02",
    "LBU_98_7",
    "LBU_05_11",
    "LBU_05_08",
    "LBU_05_03",
    "MHZ_08_04",
    "LBU_05_26",
    "LBU_98_1",
    "LBU_05_23",
    "LBU_87_2",
    "LBU_07_02",
    "LBU_05_04",
    "LBU_05_05",
    "LBU_05_12",
    "MHZ_08_03",
    "LBU_87_1",
    "LBU_96_4",
    "LBU_05_10",
    "LBU_07_01"
)

# Header cell for the new index column, matching the style already used
# by the other header cells in row 1.
$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
